$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 114.5
$ws.Range("I6").Value = 120.46667
$ws.Range("K6").Value = 361.40001
$ws.Range("M6").Value = -249.40001
$ws.Range("H9").Value = 1450.6666
$ws.Range("I9").Value = 323.1111
$ws.Range("J9").Value = 4833.3335
$ws.Range("K9").Value = 323.1111
$ws.Range("L9").Value = 4833.3335
$ws.Range("M9").Value = -154.1111
$ws.Range("N9").Value = -5171.3335
$ws.Range("H12").Value = 1205.9166
$ws.Range("J12").Value = 875
$ws.Range("L12").Value = 875
$ws.Range("N12").Value = -1215
$ws.Range("H58").Value = 1482.4
$ws.Range("J58").Value = 1798.5
$ws.Range("L58").Value = 5395.5
$ws.Range("N58").Value = -5695.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2583.8333
$ws.Range("I45").Value = 1550.9
$ws.Range("K45").Value = 1550.9
$ws.Range("M45").Value = -1173.9
$ws.Range("H61").Value = 3461.818
$ws.Range("I61").Value = 1680.3334
$ws.Range("K61").Value = 1680.3334
$ws.Range("M61").Value = -1468.3334
$ws.Range("H74").Value = 3833
$ws.Range("I74").Value = 3833
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 3833
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -2959
$ws.Range("H77").Value = 3833
$ws.Range("I77").Value = 3833
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 19165
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -14797
$ws.Range("H110").Value = 3455.5
$ws.Range("I110").Value = 2439.8
$ws.Range("K110").Value = 2439.8
$ws.Range("M110").Value = -394.8000000000002
$ws.Range("H122").Value = 2306
$ws.Range("I122").Value = 1012
$ws.Range("J122").Value = 3600
$ws.Range("K122").Value = 3036
$ws.Range("L122").Value = 10800
$ws.Range("M122").Value = -586
$ws.Range("N122").Value = -15700
$ws.Range("H136").Value = 3461.818
$ws.Range("I136").Value = 1680.3334
$ws.Range("K136").Value = 5041.0002
$ws.Range("M136").Value = -2491.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1680.75
$ws.Range("I7").Value = 1774.8334
$ws.Range("K7").Value = 1774.8334
$ws.Range("M7").Value = -1661.8334
$ws.Range("H22").Value = 822.875
$ws.Range("I22").Value = 889.125
$ws.Range("K22").Value = 889.125
$ws.Range("M22").Value = -539.125
$ws.Range("H26").Value = 7500
$ws.Range("J26").Value = 7500
$ws.Range("L26").Value = 7500
$ws.Range("N26").Value = -8074
$ws.Range("H31").Value = 3888.2856
$ws.Range("I31").Value = 3786.6667
$ws.Range("K31").Value = 3786.6667
$ws.Range("M31").Value = -3491.6667
$ws.Range("H34").Value = 3888.2856
$ws.Range("I34").Value = 3786.6667
$ws.Range("K34").Value = 3786.6667
$ws.Range("M34").Value = -3584.6667
$ws.Range("H58").Value = 1406
$ws.Range("I58").Value = 1406
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1406
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -1203
$ws.Range("H62").Value = 4335
$ws.Range("I62").Value = 3470
$ws.Range("K62").Value = 3470
$ws.Range("M62").Value = -2846
$ws.Range("H65").Value = 4335
$ws.Range("I65").Value = 3470
$ws.Range("K65").Value = 17350
$ws.Range("M65").Value = -14230
$ws.Range("H86").Value = 7747987.5
$ws.Range("I86").Value = 8715238
$ws.Range("K86").Value = 8715238
$ws.Range("M86").Value = -8714115
$ws.Range("H89").Value = 7747987.5
$ws.Range("I89").Value = 8715238
$ws.Range("K89").Value = 43576190
$ws.Range("M89").Value = -43570574
$ws.Range("H136").Value = 1406
$ws.Range("I136").Value = 1406
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4218
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -1668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 248.4
$ws.Range("J15").Value = 348
$ws.Range("L15").Value = 1044
$ws.Range("N15").Value = -1324
$ws.Range("H17").Value = 330
$ws.Range("I17").Value = 287.5
$ws.Range("J17").Value = 500
$ws.Range("K17").Value = 862.5
$ws.Range("L17").Value = 1500
$ws.Range("M17").Value = -693.5
$ws.Range("N17").Value = -1838
$ws.Range("H47").Value = 487.5
$ws.Range("I47").Value = 700
$ws.Range("K47").Value = 2100
$ws.Range("M47").Value = -1669
$ws.Range("H75").Value = 722.9091
$ws.Range("I75").Value = 950
$ws.Range("J75").Value = 593.1429000000001
$ws.Range("K75").Value = 2850
$ws.Range("L75").Value = 1779.4287
$ws.Range("M75").Value = -1852
$ws.Range("N75").Value = -3775.4287
$ws.Range("H78").Value = 722.9091
$ws.Range("I78").Value = 950
$ws.Range("J78").Value = 593.1429000000001
$ws.Range("K78").Value = 8550
$ws.Range("L78").Value = 5338.2861
$ws.Range("M78").Value = -3558
$ws.Range("N78").Value = -15322.2861
$ws.Range("H109").Value = 1220.3334
$ws.Range("I109").Value = 880.5
$ws.Range("J109").Value = 1900
$ws.Range("K109").Value = 2641.5
$ws.Range("L109").Value = 5700
$ws.Range("M109").Value = -1601.5
$ws.Range("N109").Value = -7780
$ws.Range("H112").Value = 12000
$ws.Range("I112").Value = 2000
$ws.Range("K112").Value = 6000
$ws.Range("M112").Value = -4892

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6928.5835
$ws.Range("I122").Value = 7572
$ws.Range("K122").Value = 22716
$ws.Range("M122").Value = -20266

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1296.6666
$ws.Range("I16").Value = 1293.6364
$ws.Range("K16").Value = 1293.6364
$ws.Range("M16").Value = -1123.6364
$ws.Range("H22").Value = 1316.4166
$ws.Range("I22").Value = 979.5
$ws.Range("K22").Value = 979.5
$ws.Range("M22").Value = -684.5
$ws.Range("H27").Value = 1316.4166
$ws.Range("I27").Value = 979.5
$ws.Range("K27").Value = 979.5
$ws.Range("M27").Value = -872.5
$ws.Range("H122").Value = 8428.111000000001
$ws.Range("I122").Value = 8481.625
$ws.Range("K122").Value = 25444.875
$ws.Range("M122").Value = -22994.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 933.8
$ws.Range("I113").Value = 655.1818
$ws.Range("K113").Value = 1965.5454
$ws.Range("M113").Value = 204.4546
$ws.Range("H126").Value = 2040.8182
$ws.Range("I126").Value = 1431.5
$ws.Range("K126").Value = 4294.5
$ws.Range("M126").Value = -1824.5
$ws.Range("H136").Value = 49399.273
$ws.Range("I136").Value = 65424.125
$ws.Range("K136").Value = 196272.375
$ws.Range("M136").Value = -193722.375
